$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.468507333333333
$ws.Range("H2").Value = 4.405521999999999
$ws.Range("I2").Value = 0.005118279455112885
$ws.Range("J2").Value = 0.005118279455112885
$ws.Range("M2").Value = 3.795192333333334
$ws.Range("N2").Value = 11.385577
$ws.Range("O2").Value = 0.01044213755712683
$ws.Range("P2").Value = 0.01044213755712683
$ws.Range("Q2").Value = 5.573267772910445
$ws.Range("R2").Value = 50.159409956194
$ws.Range("S2").Value = 0.00005344577812610492
$ws.Range("T2").Value = 0.00005344577812610492
$ws.Range("G3").Value = 1.468507333333333
$ws.Range("H3").Value = 4.405521999999999
$ws.Range("I3").Value = 0.005118279455112885
$ws.Range("J3").Value = 0.005118279455112885
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.6696287328350964
$ws.Range("P3").Value = 0.6696287328350964
$ws.Range("Q3").Value = 357.4000262022569
$ws.Range("R3").Value = 3216.600235820312
$ws.Range("S3").Value = 0.003427346985823149
$ws.Range("T3").Value = 0.003427346985823149
$ws.Range("G4").Value = 1.468507333333333
$ws.Range("H4").Value = 4.405521999999999
$ws.Range("I4").Value = 0.005118279455112885
$ws.Range("J4").Value = 0.005118279455112885
$ws.Range("M4").Value = 29.801371
$ws.Range("N4").Value = 89.404113
$ws.Range("O4").Value = 0.08199584844219236
$ws.Range("P4").Value = 0.08199584844219235
$ws.Range("Q4").Value = 43.76353185688733
$ws.Range("R4").Value = 393.8717867119859
$ws.Range("S4").Value = 0.000419677666486223
$ws.Range("T4").Value = 0.000419677666486223
$ws.Range("G5").Value = 1.468507333333333
$ws.Range("H5").Value = 4.405521999999999
$ws.Range("I5").Value = 0.005118279455112885
$ws.Range("J5").Value = 0.005118279455112885
$ws.Range("M5").Value = 86.47679266666667
$ws.Range("N5").Value = 259.430378
$ws.Range("O5").Value = 0.2379332811655844
$ws.Range("P5").Value = 0.2379332811655844
$ws.Range("Q5").Value = 126.9918041941462
$ws.Range("R5").Value = 1142.926237747316
$ws.Range("S5").Value = 0.001217809024677408
$ws.Range("T5").Value = 0.001217809024677408
$ws.Range("I6").Value = 0.9046276674881553
$ws.Range("J6").Value = 0.9046276674881553
$ws.Range("M6").Value = 3.795192333333334
$ws.Range("N6").Value = 11.385577
$ws.Range("O6").Value = 0.01044213755712683
$ws.Range("P6").Value = 0.01044213755712683
$ws.Range("Q6").Value = 985.0443434968099
$ws.Range("R6").Value = 8865.399091471289
$ws.Range("S6").Value = 0.009446246541894111
$ws.Range("T6").Value = 0.009446246541894111
$ws.Range("I7").Value = 0.9046276674881553
$ws.Range("J7").Value = 0.9046276674881553
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.6696287328350964
$ws.Range("P7").Value = 0.6696287328350964
$ws.Range("Q7").Value = 63168.48364748432
$ws.Range("R7").Value = 568516.3528273589
$ws.Range("S7").Value = 0.6057646786676624
$ws.Range("T7").Value = 0.6057646786676624
$ws.Range("I8").Value = 0.9046276674881553
$ws.Range("J8").Value = 0.9046276674881553
$ws.Range("M8").Value = 29.801371
$ws.Range("N8").Value = 89.404113
$ws.Range("O8").Value = 0.08199584844219236
$ws.Range("P8").Value = 0.08199584844219235
$ws.Range("Q8").Value = 7734.962909301794
$ws.Range("R8").Value = 69614.66618371614
$ws.Range("S8").Value = 0.07417571311997277
$ws.Range("T8").Value = 0.07417571311997276
$ws.Range("I9").Value = 0.9046276674881553
$ws.Range("J9").Value = 0.9046276674881553
$ws.Range("M9").Value = 86.47679266666667
$ws.Range("N9").Value = 259.430378
$ws.Range("O9").Value = 0.2379332811655844
$ws.Range("P9").Value = 0.2379332811655844
$ws.Range("Q9").Value = 22445.10105900994
$ws.Range("R9").Value = 202005.9095310894
$ws.Range("S9").Value = 0.215241029158626
$ws.Range("T9").Value = 0.215241029158626
$ws.Range("G10").Value = 0.5890733333333333
$ws.Range("H10").Value = 1.76722
$ws.Range("I10").Value = 0.002053133730501083
$ws.Range("J10").Value = 0.002053133730501083
$ws.Range("M10").Value = 3.795192333333334
$ws.Range("N10").Value = 11.385577
$ws.Range("O10").Value = 0.01044213755712683
$ws.Range("P10").Value = 0.01044213755712683
$ws.Range("Q10").Value = 2.235646598437778
$ws.Range("R10").Value = 20.12081938594
$ws.Range("S10").Value = 0.00002143910483706928
$ws.Range("T10").Value = 0.00002143910483706928
$ws.Range("G11").Value = 0.5890733333333333
$ws.Range("H11").Value = 1.76722
$ws.Range("I11").Value = 0.002053133730501083
$ws.Range("J11").Value = 0.002053133730501083
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.6696287328350964
$ws.Range("P11").Value = 0.6696287328350964
$ws.Range("Q11").Value = 143.3665464172356
$ws.Range("R11").Value = 1290.29891775512
$ws.Range("S11").Value = 0.001374837338296435
$ws.Range("T11").Value = 0.001374837338296435
$ws.Range("G12").Value = 0.5890733333333333
$ws.Range("H12").Value = 1.76722
$ws.Range("I12").Value = 0.002053133730501083
$ws.Range("J12").Value = 0.002053133730501083
$ws.Range("M12").Value = 29.801371
$ws.Range("N12").Value = 89.404113
$ws.Range("O12").Value = 0.08199584844219236
$ws.Range("P12").Value = 0.08199584844219235
$ws.Range("Q12").Value = 17.55519295287333
$ws.Range("R12").Value = 157.99673657586
$ws.Range("S12").Value = 0.0001683484421977198
$ws.Range("T12").Value = 0.0001683484421977198
$ws.Range("G13").Value = 0.5890733333333333
$ws.Range("H13").Value = 1.76722
$ws.Range("I13").Value = 0.002053133730501083
$ws.Range("J13").Value = 0.002053133730501083
$ws.Range("M13").Value = 86.47679266666667
$ws.Range("N13").Value = 259.430378
$ws.Range("O13").Value = 0.2379332811655844
$ws.Range("P13").Value = 0.2379332811655844
$ws.Range("Q13").Value = 50.94117251212889
$ws.Range("R13").Value = 458.4705526091601
$ws.Range("S13").Value = 0.0004885088451698594
$ws.Range("T13").Value = 0.0004885088451698594
$ws.Range("G14").Value = 25.306101
$ws.Range("H14").Value = 75.91830299999999
$ws.Range("I14").Value = 0.0882009193262308
$ws.Range("J14").Value = 0.0882009193262308
$ws.Range("M14").Value = 3.795192333333334
$ws.Range("N14").Value = 11.385577
$ws.Range("O14").Value = 0.01044213755712683
$ws.Range("P14").Value = 0.01044213755712683
$ws.Range("Q14").Value = 96.04152050175901
$ws.Range("R14").Value = 864.3736845158311
$ws.Range("S14").Value = 0.0009210061322695485
$ws.Range("T14").Value = 0.0009210061322695485
$ws.Range("G15").Value = 25.306101
$ws.Range("H15").Value = 75.91830299999999
$ws.Range("I15").Value = 0.0882009193262308
$ws.Range("J15").Value = 0.0882009193262308
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.6696287328350964
$ws.Range("P15").Value = 0.6696287328350964
$ws.Range("Q15").Value = 6158.907725674932
$ws.Range("R15").Value = 55430.16953107439
$ws.Range("S15").Value = 0.05906186984331449
$ws.Range("T15").Value = 0.05906186984331449
$ws.Range("G16").Value = 25.306101
$ws.Range("H16").Value = 75.91830299999999
$ws.Range("I16").Value = 0.0882009193262308
$ws.Range("J16").Value = 0.0882009193262308
$ws.Range("M16").Value = 29.801371
$ws.Range("N16").Value = 89.404113
$ws.Range("O16").Value = 0.08199584844219236
$ws.Range("P16").Value = 0.08199584844219235
$ws.Range("Q16").Value = 754.1565044644709
$ws.Range("R16").Value = 6787.408540180239
$ws.Range("S16").Value = 0.007232109213535656
$ws.Range("T16").Value = 0.007232109213535654
$ws.Range("G17").Value = 25.306101
$ws.Range("H17").Value = 75.91830299999999
$ws.Range("I17").Value = 0.0882009193262308
$ws.Range("J17").Value = 0.0882009193262308
$ws.Range("M17").Value = 86.47679266666667
$ws.Range("N17").Value = 259.430378
$ws.Range("O17").Value = 0.2379332811655844
$ws.Range("P17").Value = 0.2379332811655844
$ws.Range("Q17").Value = 2188.390449378726
$ws.Range("R17").Value = 19695.51404440853
$ws.Range("S17").Value = 0.0209859341371111
$ws.Range("T17").Value = 0.0209859341371111
